$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45911, 5913.30581487482, 5169.48725008134, 5808, 6521.068181, -1.28126599139481),
    @(3, 45912, 5913.14141862252, 4464.15489128892, 2952, 6520.665899, 88.3199738194333),
    @(4, 45913, 1177.4817126383, 2221.12242781402, 2952, 2248.96912, 14.1920764656549),
    @(5, 45914, 1060.75017668403, 2180.03479920697, 2952, 2127.948682, 12.301387688456),
    @(6, 45915, 5883.51545958832, 4971.28330018383, 2952, 6537.637226, 111.391877774813),
    @(7, 45916, 5880.24471902051, 5168.23581173641, 2952, 6722.442274, 127.434723613163),
    @(8, 45917, 6061.80298841384, 5278.82713334197, 2952, 6722.442274, 124.477767455339),
    @(9, 45918, 6061.80298841384, 5301.57173117171, 2952, 6722.442274, 125.425459031578),
    @(10, 45919, 6061.80298841384, 4552.73945347776, 2952, 6722.442274, 94.2241141276634),
    @(11, 45920, 1187.6944176276, 2210.28580401939, 2952, 2248.907556, 13.3124559329912),
    @(12, 45921, 1056.51488323291, 2141.33459009217, 2952, 2106.590983, 9.97544541080237),
    @(13, 45922, 5979.99309220113, 4949.28227899218, 2952, 6713.30322, 113.774683616294),
    @(14, 45923, 5979.99309220113, 5141.21354493961, 2952, 6713.30322, 121.771819697437),
    @(15, 45924, 5979.99309220113, 5160.24312540197, 2952, 6713.30322, 122.564718883369)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}
